# 5th test case added to BOA app
# - remove the old "expedia" sheet (no longer used)
# - keep "boaSignup" as the first sheet
# - add a new "loginNegativeTest" sheet with negative-login test data

$wb = $excel.ActiveWorkbook

# Drop the now-unused "expedia" sheet.
$expedia = $wb.Worksheets.Item("expedia")
$expedia.Delete() | Out-Null

$boaSignup = $wb.Worksheets.Item("boaSignup")

# Add the new sheet right after boaSignup; it becomes the active/selected tab.
$loginNegativeTest = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $boaSignup)
$loginNegativeTest.Name = "loginNegativeTest"

# Fill data column-by-column so shared strings are interned in the same
# order the original author would have typed them in.
$loginNegativeTest.Range("A1").Value = "OnlineId"
$loginNegativeTest.Range("B1").Value = "Password"

$loginNegativeTest.Range("A2").Value = "ABCDE"
$loginNegativeTest.Range("A3").Value = "FGHIJK"
$loginNegativeTest.Range("A4").Value = "LMNOPQ"

$loginNegativeTest.Range("B2").Value = "abcd123#"
$loginNegativeTest.Range("B3").Value = "FHG1234$"
$loginNegativeTest.Range("B4").Value = "pqrst123#"

$loginNegativeTest.Range("C1").Value = "ErrContains"
$loginNegativeTest.Range("C2").Value = "does not match"
$loginNegativeTest.Range("C3").Value = "does not match"
$loginNegativeTest.Range("C4").Value = "does not match"

$loginNegativeTest.Range("D12").Select() | Out-Null
